$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C24').Value = 'Openness'
$ws.Range('A35').Value = 'Urgency/Priority'
$ws.Range('C89').Value = 'Existing Policies'
$ws.Range('C124').ClearContents()
$ws.Range('C128').Value = 'Stakeholder Opposition'
$ws.Range('C138').Value = 'Existing Policies'
$ws.Range('B145').Value = 'Existing Policies'
$ws.Range('C151').ClearContents()
$ws.Range('C164').Value = 'Human'
$ws.Range('B169').Value = 'Human'
$ws.Range('B173').Value = 'Info Share/Use'
$ws.Range('C175').Value = 'Reg Leader'
$ws.Range('C181').Value = 'Competition'
$ws.Range('C194').Value = 'Reg Leader'
$ws.Range('B196').Value = 'Urgency/Priority'
$ws.Range('B209').Value = 'Financial'
$ws.Range('C210').Value = 'Info Share/Use'
$ws.Range('C213').Value = 'Urgency/Priority'
$ws.Range('A224').Value = 'Reg Leader'
$ws.Range('B233').Value = 'Financial'
$ws.Range('A242').Value = 'Stakeholder Engage/Influence'
$ws.Range('C262').Value = 'Info Share/Use'
$ws.Range('B268').Value = 'Collab Experience'
$ws.Range('B280').Value = 'Info Share/Use'
$ws.Range('B290').Value = 'Human'
$ws.Range('B292').ClearContents()
$ws.Range('C294').Value = 'Overall Plan'
$ws.Range('C298').Value = 'Urgency/Priority'
$ws.Range('C300').Value = 'Stakeholder Engage/Influence'
$ws.Range('C301').Value = 'Competition'
$ws.Range('B303').Value = 'Outside Scope'
$ws.Range('A324').Value = 'Narrow Grey Focus'
$ws.Range('C331').Value = 'Reg Leader'
$ws.Range('A334').Value = 'Urgency/Priority'
$ws.Range('B340').Value = 'Competition'
$ws.Range('C354').Value = 'Narrow Grey Focus'
$ws.Range('C356').ClearContents()
$ws.Range('B361').ClearContents()
$ws.Range('A382').Value = 'Outside Scope'
$ws.Range('C387').Value = 'Narrow Grey Focus'
$ws.Range('A388').Value = 'Outside Scope'
$ws.Range('C396').ClearContents()
$ws.Range('A399').Value = 'Political Leader'
$ws.Range('A401').Value = 'Urgency/Priority'
$ws.Range('A419').Value = 'Reg Leader'
$ws.Range('C427').ClearContents()
$ws.Range('C429').Value = 'Info Share/Use'
$ws.Range('B436').Value = 'Human'
$ws.Range('C444').Value = 'Competition'
$ws.Range('C448').Value = 'Urgency/Priority'
$ws.Range('C460').Value = 'Human'
$ws.Range('C466').Value = 'Human'
$ws.Range('C467').Value = 'Urgency/Priority'
$ws.Range('C468').Value = 'Competition'
$ws.Range('A472').Value = 'Partner Capacity'
$ws.Range('B477').ClearContents()
$ws.Range('C499').Value = 'Reg Leader'
$ws.Range('B507').Value = 'Human'
$ws.Range('B523').Value = 'Financial'
$ws.Range('C523').Value = 'Partner Capacity'
$ws.Range('A523').Value = 'Human'
$ws.Range('C547').Value = 'Urgency/Priority'
$ws.Range('A552').Value = 'Urgency/Priority'
$ws.Range('C561').Value = 'Org Leader'
$ws.Range('B564').Value = 'Stakeholder Engage/Influence'
$ws.Range('A568').Value = 'Outside Scope'
$ws.Range('C569').Value = 'Stakeholder Engage/Influence'
$ws.Range('C579').Value = 'Urgency/Priority'
$ws.Range('C584').Value = 'Human'
$ws.Range('A587').Value = 'Outside Scope'
$ws.Range('B590').Value = 'Human'
$ws.Range('A590').Value = 'Outside Scope'
$ws.Range('C594').Value = 'Urgency/Priority'
$ws.Range('A595').Value = 'Outside Scope'
$ws.Range('C615').Value = 'Openness'
$ws.Range('C616').Value = 'Stakeholder Engage/Influence'
$ws.Range('B622').Value = 'Human'
$ws.Range('B623').Value = 'Partner Capacity'
$ws.Range('B626').Value = 'Stakeholder Engage/Influence'
$ws.Range('A628').Value = 'Outside Scope'
$ws.Range('A632').Value = 'Outside Scope'
$ws.Range('B635').Value = 'Stakeholder Engage/Influence'
$ws.Range('C641').Value = 'Stakeholder Engage/Influence'
$ws.Range('B651').Value = 'Financial'
$ws.Range('B659').Value = 'Collab Experience'
$ws.Range('C666').Value = 'SLR Uncertainty'
$ws.Range('C677').Value = 'Info Share/Use'
$ws.Range('B690').Value = 'Human'
$ws.Range('C700').Value = 'Stakeholder Engage/Influence'
$ws.Range('C705').Value = 'Human'
$ws.Range('C709').Value = 'Info Share/Use'
$ws.Range('A723').Value = 'Outside Scope'
$ws.Range('C730').Value = 'Urgency/Priority'
$ws.Range('C734').Value = 'Openness'
$ws.Range('C741').Value = 'Collab Experience'
$ws.Range('A747').Value = 'Stakeholder Engage/Influence'
$ws.Range('C748').Value = 'Narrow Grey Focus'
$ws.Range('C768').Value = 'SLR Uncertainty'
$ws.Range('B769').ClearContents()
$ws.Range('B779').Value = 'Urgency/Priority'
$ws.Range('C792').Value = 'Info Share/Use'
$ws.Range('A794').Value = 'Urgency/Priority'
$ws.Range('A795').Value = 'Urgency/Priority'
$ws.Range('C850').Value = 'Financial'

$ws.Range('A259').ClearContents()
$ws.Range('A485').ClearContents()
$ws.Range('A725').ClearContents()
$ws.Range('A739').ClearContents()

$ws.Range('A747').Interior.ColorIndex = -4142

[void]$ws.Range('H5').Select()